# TER2021-DOW-template.docx edits
$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function ReplaceText($search, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $search"
    } else {
        Write-Host "OK: $search"
    }
}

# 1. "Chafoun" -> "Chalfoun" everywhere (Etudiant(s) list + Organisation du travail list)
ReplaceText "Chafoun" "Chalfoun"

# 2. "les bases données qu'il souhaite analyser" -> "les bases de données qu'il souhaite analyser"
ReplaceText "les bases données qu’il souhaite analyser" "les bases de données qu’il souhaite analyser"

# 3. " faites. En effet,  nous devrons vérifier" -> " faites. Nous devrons vérifier"
ReplaceText " faites. En effet,  nous devrons vérifier" " faites. Nous devrons vérifier"

# 4. "Enfin, l'utilisateur pourrait télécharger" -> "Enfin, l'utilisateur pourra télécharger"
ReplaceText "Enfin, l’utilisateur pourrait télécharger" "Enfin, l’utilisateur pourra télécharger"

# 5. "des développeurs,  data " -> "des développeurs, data "
ReplaceText "des développeurs,  data " "des développeurs, data "

# 6. "Quel utilité a le projet<nbsp>? " -> "Quelle utilité a le projet<nbsp>?]"
$s6 = "Quel utilité a le projet" + $nbsp + "? "
$r6 = "Quelle utilité a le projet" + $nbsp + "?]"
ReplaceText $s6 $r6

# 7. "Assurer la sécurité => JEREMY" -> full sentence
ReplaceText "Assurer la sécurité => JEREMY" "Assurer une protection contre les cyberattaques (notamment XSS, Session Management, SQL attack etc…)"

# 8. "d'une base données, l'algorithme" -> "d'une base de données, l'algorithme"
ReplaceText "d’une base données, l’algorithme" "d’une base de données, l’algorithme"

# 9. "Le drag and drop est une porte d'entrée pour tout type d'attaque. " -> add "de base de données" + new sentence
ReplaceText "Le drag and drop est une porte d’entrée pour tout type d’attaque. " "Le drag and drop de base de données est une porte d’entrée pour tout type d’attaque. Une attention particulière devra être portée sur la sécurité, à la fois pour le serveur de l’application mais également pour les futurs utilisateurs qui feront confiance au site en uploadant leurs bases de données. "

# 10. Append sentence about plafond de requête to "La surcharge du serveur..." paragraph
ReplaceText "La surcharge du serveur qui peut entraîner une inutilisation de l’application. " "La surcharge du serveur qui peut entraîner une inutilisation de l’application. En fonction de la puissance du serveur nous devront établir un plafond de requête par utilisateur."

# 11. Insert a new bullet paragraph before "Une base données trop compliquée..." and fix "base données" -> "base de données"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Une base données trop compliquée")) {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphBefore() | Out-Null
$newPara = $target.Previous()
$newPara.Range.Text = "Le volume des bases de données. Une taille maximale de fichier devra être fixée afin d’en un premier temps d’éviter les temps trop longs de téléchargement et réduire le nombre de calcul. Trouver le système d’envois le plus adapté (compression, format, …)"

ReplaceText "Une base données trop compliquée à analyser pour un modèle de Machine Learning. " "Une base de données trop compliquée à analyser pour un modèle de Machine Learning. "

# 12. "il pourra le faire<nbsp>! En effet, imaginons l'étudiant ne comprend pas bien l'intérêt du " -> "... ! Imaginons que l'étudiant ne comprenne pas bien l'intérêt du "
$s12 = "il pourra le faire" + $nbsp + "! En effet, imaginons l’étudiant ne comprend pas bien l’intérêt du "
$r12 = "il pourra le faire" + $nbsp + "! Imaginons que l’étudiant ne comprenne pas bien l’intérêt du "
ReplaceText $s12 $r12

# 13. "- Ralph El Chalfoun s'occupe de la partie développement de l'application web" -> add trailing "."
ReplaceText "- Ralph El Chalfoun s’occupe de la partie développement de l’application web" "- Ralph El Chalfoun s’occupe de la partie développement de l’application web."

Write-Host "All replacements attempted."
